# Auto-generated edit script: update crypto price/volume table (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Text" number format on price cells whose new values look numeric,
# so COM stores the literal digit string instead of auto-converting to a Number.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values
$ws.Range("D2").Value = "27.056.97"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.864.97"
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "306.37"
$ws.Range("E5").Value = "  -1.91%  "
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "0.5134"
$ws.Range("E7").Value = "  -1.83%  "
$ws.Range("D8").Value = "0.3757"
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "0.07165"
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "0.8901"
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("D11").Value = "20.70"
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("D12").Value = "0.07593"
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "1.849.48"
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").Value = "5.311"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "89.44"
$ws.Range("E15").Value = "  -3.02%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("D17").Value = "0.000008462"
$ws.Range("E17").Value = "  -2.73%  "
$ws.Range("D18").Value = "14.12"
$ws.Range("E18").Value = "  -2.69%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("D20").Value = "27.058.93"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("D21").Value = "5.041"
$ws.Range("E21").Value = "  -2.00%  "
$ws.Range("D22").Value = "2.091.92"
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D24").Value = "6.452"
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").Value = "147.48"
$ws.Range("E26").Value = "  -3.68%  "
$ws.Range("D27").Value = "18.00"
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").Value = "2.114"
$ws.Range("E28").Value = "  -2.61%  "
$ws.Range("D29").Value = "112.77"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "4.660"
$ws.Range("E30").Value = "  -4.19%  "
$ws.Range("D31").Value = "4.708"
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("D32").Value = "0.09114"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "0.05141"
$ws.Range("E33").Value = "  -2.79%  "
$ws.Range("D34").Value = "3.063"
$ws.Range("E34").Value = "  -3.45%  "
$ws.Range("E35").Value = "  -6.13%  "
$ws.Range("D36").Value = "0.7270"
$ws.Range("E36").Value = "  -6.11%  "
$ws.Range("D37").Value = "0.02041"
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("D39").Value = "2.473"
$ws.Range("E39").Value = "  -6.23%  "
$ws.Range("D40").Value = "1.076"
$ws.Range("E40").Value = "  -1.58%  "
$ws.Range("D41").Value = "0.5337"
$ws.Range("E41").Value = "  -3.23%  "
$ws.Range("D42").Value = "6.568"
$ws.Range("E42").Value = "  -1.60%  "
$ws.Range("D43").Value = "117.24"
$ws.Range("E43").Value = "  +2.32%  "
$ws.Range("D44").Value = "8.278"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").Value = "0.1473"
$ws.Range("E45").Value = "  -2.54%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.000"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4633"
$ws.Range("E47").Value = "  -3.72%  "
$ws.Range("D48").Value = "10.01"
$ws.Range("E48").Value = "  -4.01%  "
$ws.Range("D49").Value = "1.577"
$ws.Range("E49").Value = "  -2.41%  "
$ws.Range("D50").Value = "36.63"
$ws.Range("E50").Value = "  -0.69%  "
$ws.Range("D51").Value = "63.84"
$ws.Range("E51").Value = "  -4.44%  "

Write-Host "Updated cryptos list"
